$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.07024999999999999
$ws.Range("H2").Value = 0.21075
$ws.Range("I2").Value = 0.005236595731231519
$ws.Range("J2").Value = 0.005236595731231519
$ws.Range("M2").Value = 1.599392
$ws.Range("N2").Value = 4.798176
$ws.Range("O2").Value = 0.03952976301548796
$ws.Range("P2").Value = 0.03952976301548796
$ws.Range("Q2").Value = 0.112357288
$ws.Range("R2").Value = 1.011215592
$ws.Range("S2").Value = 0.0002070013882634978
$ws.Range("T2").Value = 0.0002070013882634978

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.07024999999999999
$ws.Range("H3").Value = 0.21075
$ws.Range("I3").Value = 0.005236595731231519
$ws.Range("J3").Value = 0.005236595731231519
$ws.Range("O3").Value = 0.4638329693976876
$ws.Range("P3").Value = 0.4638329693976876
$ws.Range("Q3").Value = 1.318374069333333
$ws.Range("R3").Value = 11.865366624
$ws.Range("S3").Value = 0.002428905747552371
$ws.Range("T3").Value = 0.002428905747552371

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.07024999999999999
$ws.Range("H4").Value = 0.21075
$ws.Range("I4").Value = 0.005236595731231519
$ws.Range("J4").Value = 0.005236595731231519
$ws.Range("M4").Value = 20.09416733333333
$ws.Range("N4").Value = 60.28250199999999
$ws.Range("O4").Value = 0.4966372675868244
$ws.Range("P4").Value = 0.4966372675868245
$ws.Range("Q4").Value = 1.411615255166666
$ws.Range("R4").Value = 12.7045372965
$ws.Range("S4").Value = 0.00260068859541565
$ws.Range("T4").Value = 0.002600688595415651

$ws.Range("I5").Value = 0.4287876899474159
$ws.Range("J5").Value = 0.4287876899474159
$ws.Range("M5").Value = 1.599392
$ws.Range("N5").Value = 4.798176
$ws.Range("O5").Value = 0.03952976301548796
$ws.Range("P5").Value = 0.03952976301548796
$ws.Range("Q5").Value = 9.200141550538666
$ws.Range("R5").Value = 82.801273954848
$ws.Range("S5").Value = 0.01694987576757988
$ws.Range("T5").Value = 0.01694987576757988

$ws.Range("I6").Value = 0.4287876899474159
$ws.Range("J6").Value = 0.4287876899474159
$ws.Range("O6").Value = 0.4638329693976876
$ws.Range("P6").Value = 0.4638329693976876
$ws.Range("S6").Value = 0.1988858674694849
$ws.Range("T6").Value = 0.1988858674694849

$ws.Range("I7").Value = 0.4287876899474159
$ws.Range("J7").Value = 0.4287876899474159
$ws.Range("M7").Value = 20.09416733333333
$ws.Range("N7").Value = 60.28250199999999
$ws.Range("O7").Value = 0.4966372675868244
$ws.Range("P7").Value = 0.4966372675868245
$ws.Range("Q7").Value = 115.5871630012384
$ws.Range("R7").Value = 1040.284467011146
$ws.Range("S7").Value = 0.2129519467103511
$ws.Range("T7").Value = 0.2129519467103511

$ws.Range("G8").Value = 7.592679666666666
$ws.Range("H8").Value = 22.778039
$ws.Range("I8").Value = 0.5659757143213526
$ws.Range("J8").Value = 0.5659757143213525
$ws.Range("M8").Value = 1.599392
$ws.Range("N8").Value = 4.798176
$ws.Range("O8").Value = 0.03952976301548796
$ws.Range("P8").Value = 0.03952976301548796
$ws.Range("Q8").Value = 12.14367111742933
$ws.Range("R8").Value = 109.293040056864
$ws.Range("S8").Value = 0.02237288585964458
$ws.Range("T8").Value = 0.02237288585964458

$ws.Range("G9").Value = 7.592679666666666
$ws.Range("H9").Value = 22.778039
$ws.Range("I9").Value = 0.5659757143213526
$ws.Range("J9").Value = 0.5659757143213525
$ws.Range("O9").Value = 0.4638329693976876
$ws.Range("P9").Value = 0.4638329693976876
$ws.Range("Q9").Value = 142.4909891713565
$ws.Range("R9").Value = 1282.418902542208
$ws.Range("S9").Value = 0.2625181961806503
$ws.Range("T9").Value = 0.2625181961806503

$ws.Range("G10").Value = 7.592679666666666
$ws.Range("H10").Value = 22.778039
$ws.Range("I10").Value = 0.5659757143213526
$ws.Range("J10").Value = 0.5659757143213525
$ws.Range("M10").Value = 20.09416733333333
$ws.Range("N10").Value = 60.28250199999999
$ws.Range("O10").Value = 0.4966372675868244
$ws.Range("P10").Value = 0.4966372675868245
$ws.Range("Q10").Value = 152.5685757303975
$ws.Range("R10").Value = 1373.117181573578
$ws.Range("S10").Value = 0.2810846322810577
$ws.Range("T10").Value = 0.2810846322810577
